$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - first data block
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 192
$ws1.Range("F8").Value = 276
$ws1.Range("F10").Value = 1029
$ws1.Range("F15").Value = 13093
$ws1.Range("F16").Value = 159
$ws1.Range("F18").Value = 14
$ws1.Range("F19").Value = 5377
$ws1.Range("F20").Value = 5551
$ws1.Range("F21").Value = 9

# Sheet "全部类型" (All Types) - second data block
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 192
$ws4.Range("F24").Value = 276
$ws4.Range("F32").Value = 1029
$ws4.Range("F37").Value = 13093
$ws4.Range("F38").Value = 159
$ws4.Range("F41").Value = 14
$ws4.Range("F42").Value = 5377
$ws4.Range("F43").Value = 5551
$ws4.Range("F44").Value = 9

$wb.Save()
